$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------

# Row 3 used to hold "Minimal-Variance matching", which is removed entirely;
# the row below it ("Elastic action comparison with freedom degree") shifts
# up to become the new row 3, keeping its own formatting (e.g. its 30pt
# row height).
$ws.Rows.Item(3).Delete()

# Header row: G1/H1/I1 used to contain "Right hand punch"/"Left hand punch"/"Jump".
# Now only G1 remains used, holding the new "Taekwondo position" entry; H1/I1
# become blank (but keep their non-bold header style).
$ws.Cells.Item(1, 7).Value = "Taekwondo position"
$ws.Cells.Item(1, 8).Value = ""
$ws.Cells.Item(1, 9).Value = ""

# --- New numeric data rows --------------------------------------------

$ws.Cells.Item(2, 2).Value = 0.866
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 0.9166
$ws.Cells.Item(2, 6).Value = 0.7143

$ws.Cells.Item(3, 2).Value = 0.8
$ws.Cells.Item(3, 3).Value = 0.75
$ws.Cells.Item(3, 4).Value = 0.8182
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.7557

# Number formats: column C/D of row 2 use whole-percent, everything else
# that holds a ratio uses percent with two decimals.
$ws.Cells.Item(2, 2).NumberFormat = "0.00%"
$ws.Cells.Item(2, 3).NumberFormat = "0%"
$ws.Cells.Item(2, 4).NumberFormat = "0%"
$ws.Cells.Item(2, 5).NumberFormat = "0.00%"
$ws.Cells.Item(2, 6).NumberFormat = "0.00%"

$ws.Cells.Item(3, 2).NumberFormat = "0.00%"
$ws.Cells.Item(3, 3).NumberFormat = "0.00%"
$ws.Cells.Item(3, 4).NumberFormat = "0.00%"
$ws.Cells.Item(3, 5).NumberFormat = "0.00%"
$ws.Cells.Item(3, 6).NumberFormat = "0.00%"

# --- Bold styling for row/column headers ------------------------------

# A1 has no text but still picks up the bold header formatting.
$ws.Cells.Item(1, 1).Font.Bold = $true

$ws.Range("B1:G1").Font.Bold = $true
$ws.Range("A2:A3").Font.Bold = $true

# --- Column width & selection ------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 6.71

$ws.Range("H1").Select() | Out-Null
